$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "key" sub-field for each service entry: the search logic used to
# reuse the category keys (cat11, cat12, ...) for the nested services, which
# broke lookups. Rename them to their own dedicated "srv" keys so they no
# longer collide with the category keys (cat1, cat2, ...).

$ws.Range("F2").Value = "srv11"
$ws.Range("K2").Value = "srv12"
$ws.Range("P2").Value = "srv13"
$ws.Range("U2").Value = "srv14"
$ws.Range("Z2").Value = "srv15"
$ws.Range("AE2").Value = "srv16"

$ws.Range("F3").Value = "srv21"
$ws.Range("K3").Value = "srv22"
$ws.Range("P3").Value = "srv23"
$ws.Range("U3").Value = "srv24"
$ws.Range("Z3").Value = "srv25"
$ws.Range("AE3").Value = "srv26"
$ws.Range("AJ3").Value = "srv27"
$ws.Range("AO3").Value = "srv28"
$ws.Range("AT3").Value = "srv29"

$ws.Range("F4").Value = "srv31"
$ws.Range("K4").Value = "srv32"
$ws.Range("P4").Value = "srv33"
$ws.Range("U4").Value = "srv34"
$ws.Range("Z4").Value = "srv35"
$ws.Range("AE4").Value = "srv36"

$ws.Range("F5").Value = "srv41"
$ws.Range("K5").Value = "srv42"
$ws.Range("P5").Value = "srv43"
$ws.Range("U5").Value = "srv44"

# Move the active selection to A6, matching where the author left off.
$ws.Range("A6").Select()
